# Notifications.xlsx - "Added test cases in to Notification module"
# Appends three new rows (Notifications0017/0018/0019) to the "Test Cases"
# worksheet, copying the formatting of the last existing data row (17) and
# then filling in the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Copy the formatting (style) of the last populated row onto the new rows
# before writing values into them, so the new cells inherit the same cell
# style (fill/border) as the rest of the table.
$ws.Range("A17:E17").Copy()
$ws.Range("A18:E20").PasteSpecial(-4122)

# Row 18 - Notifications0017
$ws.Range("A18").Value = "Notifications0017"
$ws.Range("B18").Value = "OPQA-1601"
$ws.Range("C18").Value = "Verify that user is able to navigate record view page by clicking article title from Recommended articles section on Home page"
$ws.Range("D18").Value = "Y"
$ws.Range("E18").Value = "PASS"

# Row 19 - Notifications0018
$ws.Range("A19").Value = "Notifications0018"
$ws.Range("B19").Value = "OPQA-1602"
$ws.Range("C19").Value = "Verify that user is able to watch article from Recommended articles section on Home page."
$ws.Range("D19").Value = "Y"
$ws.Range("E19").Value = "PASS"

# Row 20 - Notifications0019
$ws.Range("A20").Value = "Notifications0019"
$ws.Range("B20").Value = "OPQA-1600"
$ws.Range("C20").Value = 'Verify that user ia able to publish post by clicking "Publish a post of your own" link Feature post section on Home page.'
$ws.Range("D20").Value = "Y"
$ws.Range("E20").Value = "PASS"

# Keep the sheet selection/view in sync with the newly extended data range.
[void]$ws.Range("D2:D20").Select()
